$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$attachBase = "https://github.com/Oleksandr-Mnk/Test-documentation/tree/main/Bug%20reports/Attachments%20to%20bug%20reports"

# Match the existing "hyperlink cell" formatting already used by J2:J9
# (left/centered alignment, no wrap) before the hyperlinks are wired up,
# so the newly linked cells end up visually identical to the others.
$rng = $ws.Range("J10:J12")
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108
$rng.WrapText = $false

# Set the display text for the newly-linked cells so that the individual
# hyperlinks below can be added without re-specifying "TextToDisplay"
# (keeps them referencing the existing shared string "Link to attachments").
$ws.Range("J10").Value = "Link to attachments"
$ws.Range("J11").Value = "Link to attachments"
$ws.Range("J12").Value = "Link to attachments"

# Range hyperlink covering J10:J12 pointing at the attachments folder.
$ws.Hyperlinks.Add($ws.Range("J10:J12"), $attachBase, "", "", "Link to attachments")

# Individual hyperlinks for each bug report's attachment folder.
$ws.Hyperlinks.Add($ws.Range("J10"), "$attachBase/BR9", "", "", "")
$ws.Hyperlinks.Add($ws.Range("J11"), "$attachBase/BR10", "", "", "")
$ws.Hyperlinks.Add($ws.Range("J12"), "$attachBase/BR11", "", "", "")

# Restore the view to reflect scrolling down and selecting B15.
$ws.Range("B15").Select() | Out-Null
